$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 844.41174
$ws.Range("I15").Value = 844.41174
$ws.Range("K15").Value = 2533.23522
$ws.Range("M15").Value = -2364.23522
$ws.Range("H17").Value = 916.7857
$ws.Range("J17").Value = 914.0755
$ws.Range("L17").Value = 2742.2265
$ws.Range("N17").Value = -3078.2265
$ws.Range("H62").Value = 18360.385
$ws.Range("I62").Value = 17698.092
$ws.Range("K62").Value = 17698.092
$ws.Range("M62").Value = -17074.092
$ws.Range("H65").Value = 18360.385
$ws.Range("I65").Value = 17698.092
$ws.Range("K65").Value = 88490.46000000001
$ws.Range("M65").Value = -85370.46000000001
$ws.Range("H88").Value = 4462.2354
$ws.Range("I88").Value = 1987.375
$ws.Range("K88").Value = 1987.375
$ws.Range("M88").Value = -1581.375
$ws.Range("H91").Value = 4462.2354
$ws.Range("I91").Value = 1987.375
$ws.Range("K91").Value = 1987.375
$ws.Range("M91").Value = -583.375
$ws.Range("H106").Value = 5162.7856
$ws.Range("I106").Value = 5028.5
$ws.Range("J106").Value = 5498.5
$ws.Range("K106").Value = 5028.5
$ws.Range("L106").Value = 5498.5
$ws.Range("M106").Value = -4397.5
$ws.Range("N106").Value = -6760.5
$ws.Range("H107").Value = 400.2414
$ws.Range("I107").Value = 395.65216
$ws.Range("J107").Value = 417.83334
$ws.Range("K107").Value = 395.65216
$ws.Range("L107").Value = 417.83334
$ws.Range("M107").Value = 1524.34784
$ws.Range("N107").Value = -4257.83334
$ws.Range("H111").Value = 4081.7778
$ws.Range("I111").Value = 3100.5
$ws.Range("K111").Value = 9301.5
$ws.Range("M111").Value = -6234.5
$ws.Range("H112").Value = 1137.2439
$ws.Range("J112").Value = 1153.3684
$ws.Range("L112").Value = 3460.1052
$ws.Range("N112").Value = -5676.1052
$ws.Range("H132").Value = 35313.52
$ws.Range("J132").Value = 6502.5386
$ws.Range("L132").Value = 19507.6158
$ws.Range("N132").Value = -24567.6158
$ws.Range("H138").Value = 7162.3193
$ws.Range("J138").Value = 2790.2068
$ws.Range("L138").Value = 8370.6204
$ws.Range("N138").Value = -18650.6204
$ws.Range("H141").Value = 1322.6
$ws.Range("I141").Value = 1322.6
$ws.Range("K141").Value = 3967.8
$ws.Range("M141").Value = 1212.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6411093
$ws.Range("I32").Value = 6579782.5
$ws.Range("K32").Value = 6579782.5
$ws.Range("M32").Value = -6579495.5
$ws.Range("H45").Value = 1982.4642
$ws.Range("I45").Value = 1943.84
$ws.Range("K45").Value = 1943.84
$ws.Range("M45").Value = -1566.84
$ws.Range("H61").Value = 1077665.1
$ws.Range("I61").Value = 1668265.6
$ws.Range("K61").Value = 1668265.6
$ws.Range("M61").Value = -1668053.6
$ws.Range("H63").Value = 24396
$ws.Range("I63").Value = 5493.25
$ws.Range("K63").Value = 5493.25
$ws.Range("M63").Value = -4807.25
$ws.Range("H66").Value = 24396
$ws.Range("I66").Value = 5493.25
$ws.Range("K66").Value = 27466.25
$ws.Range("M66").Value = -24034.25
$ws.Range("H110").Value = 2008.8
$ws.Range("I110").Value = 2633.3333
$ws.Range("K110").Value = 2633.3333
$ws.Range("M110").Value = -588.3332999999998
$ws.Range("H136").Value = 1077665.1
$ws.Range("I136").Value = 1668265.6
$ws.Range("K136").Value = 5004796.800000001
$ws.Range("M136").Value = -5002246.800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1442.8462
$ws.Range("I86").Value = 1299.3334
$ws.Range("J86").Value = 1765.75
$ws.Range("K86").Value = 1299.3334
$ws.Range("L86").Value = 1765.75
$ws.Range("M86").Value = -176.3334
$ws.Range("N86").Value = -4011.75
$ws.Range("H89").Value = 1442.8462
$ws.Range("I89").Value = 1299.3334
$ws.Range("J89").Value = 1765.75
$ws.Range("K89").Value = 6496.666999999999
$ws.Range("L89").Value = 8828.75
$ws.Range("M89").Value = -880.6669999999995
$ws.Range("N89").Value = -20060.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 126996.93
$ws.Range("I31").Value = 261795.25
$ws.Range("K31").Value = 261795.25
$ws.Range("M31").Value = -261500.25
$ws.Range("H34").Value = 126996.93
$ws.Range("I34").Value = 261795.25
$ws.Range("K34").Value = 261795.25
$ws.Range("M34").Value = -261593.25
$ws.Range("H58").Value = 414959.56
$ws.Range("I58").Value = 618644.6
$ws.Range("J58").Value = 7589.5
$ws.Range("K58").Value = 618644.6
$ws.Range("L58").Value = 7589.5
$ws.Range("M58").Value = -618441.6
$ws.Range("N58").Value = -7995.5
$ws.Range("H107").Value = 1199.1666
$ws.Range("I107").Value = 1249
$ws.Range("J107").Value = 1174.25
$ws.Range("K107").Value = 1249
$ws.Range("L107").Value = 1174.25
$ws.Range("M107").Value = 671
$ws.Range("N107").Value = -5014.25
$ws.Range("H136").Value = 414959.56
$ws.Range("I136").Value = 618644.6
$ws.Range("J136").Value = 7589.5
$ws.Range("K136").Value = 1855933.8
$ws.Range("L136").Value = 22768.5
$ws.Range("M136").Value = -1853383.8
$ws.Range("N136").Value = -27868.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 6930.5713
$ws.Range("I99").Value = 1951
$ws.Range("J99").Value = 8922.4
$ws.Range("K99").Value = 5853
$ws.Range("L99").Value = 26767.2
$ws.Range("M99").Value = -3607
$ws.Range("N99").Value = -31259.2
$ws.Range("H141").Value = 1780.7
$ws.Range("I141").Value = 1578.5555
$ws.Range("K141").Value = 4735.666499999999
$ws.Range("M141").Value = 444.3335000000006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1152.2941
$ws.Range("I102").Value = 1265.75
$ws.Range("J102").Value = 880
$ws.Range("K102").Value = 1265.75
$ws.Range("L102").Value = 880
$ws.Range("M102").Value = 356.25
$ws.Range("N102").Value = -4124
$ws.Range("H113").Value = 3998.6
$ws.Range("I113").Value = 3998
$ws.Range("K113").Value = 3998
$ws.Range("M113").Value = -1828

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 450.52
$ws.Range("I16").Value = 439.26666
$ws.Range("J16").Value = 467.4
$ws.Range("K16").Value = 439.26666
$ws.Range("L16").Value = 467.4
$ws.Range("M16").Value = -269.26666
$ws.Range("N16").Value = -807.4
$ws.Range("H22").Value = 507.85715
$ws.Range("I22").Value = 435.5
$ws.Range("J22").Value = 604.3333
$ws.Range("K22").Value = 435.5
$ws.Range("L22").Value = 604.3333
$ws.Range("M22").Value = -140.5
$ws.Range("N22").Value = -1194.3333
$ws.Range("H27").Value = 507.85715
$ws.Range("I27").Value = 435.5
$ws.Range("J27").Value = 604.3333
$ws.Range("K27").Value = 435.5
$ws.Range("L27").Value = 604.3333
$ws.Range("M27").Value = -328.5
$ws.Range("N27").Value = -818.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 68005
$ws.Range("J46").Value = 68005
$ws.Range("L46").Value = 68005
$ws.Range("N46").Value = -68467
$ws.Range("H107").Value = 2035.5652
$ws.Range("I107").Value = 1169.2941
$ws.Range("J107").Value = 4490
$ws.Range("K107").Value = 3507.8823
$ws.Range("L107").Value = 13470
$ws.Range("M107").Value = -1587.8823
$ws.Range("N107").Value = -17310
$ws.Range("H113").Value = 1626.2222
$ws.Range("I113").Value = 477.25
$ws.Range("J113").Value = 2545.4
$ws.Range("K113").Value = 1431.75
$ws.Range("L113").Value = 7636.200000000001
$ws.Range("M113").Value = 738.25
$ws.Range("N113").Value = -11976.2
$ws.Range("H134").Value = 68005
$ws.Range("J134").Value = 68005
$ws.Range("L134").Value = 204015
$ws.Range("N134").Value = -209085

Write-Host "Updated 194 cells"